# Add a new IPPIS_NO column at the start of the sheet (insert before column A),
# shifting the existing header columns one to the right.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(1).Insert()

$ws.Range("A1").Value = "IPPIS_NO"

# Move the selection as recorded in the saved workbook
$ws.Range("A8").Select()
